$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5514945387840271
$ws.Range("B1").Value = 0.693696916103363
$ws.Range("C1").Value = 1.003373146057129
$ws.Range("D1").Value = 4.112203121185303
$ws.Range("E1").Value = 3.980604648590088
